$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-04-11 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-04-12 Saturday", 2)

# Update the multiplication-table answers. The table has 20 rows, but only
# rows 1, 5, 10, 15 and 20 contain the 5 answer cells (the other rows are
# blank spacer rows), so we address each answer cell explicitly by its
# row/column position to avoid ambiguity from duplicate "before" values
# (e.g. "40×73=2920" appears twice but maps to two different results).
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="31×88=2728"},
    @{Row=1;  Col=2; Text="24×75=1800"},
    @{Row=1;  Col=3; Text="91×71=6461"},
    @{Row=1;  Col=4; Text="15×87=1305"},
    @{Row=1;  Col=5; Text="96×19=1824"},

    @{Row=5;  Col=1; Text="64×58=3712"},
    @{Row=5;  Col=2; Text="53×91=4823"},
    @{Row=5;  Col=3; Text="58×31=1798"},
    @{Row=5;  Col=4; Text="81×96=7776"},
    @{Row=5;  Col=5; Text="54×53=2862"},

    @{Row=10; Col=1; Text="65×95=6175"},
    @{Row=10; Col=2; Text="51×61=3111"},
    @{Row=10; Col=3; Text="55×37=2035"},
    @{Row=10; Col=4; Text="96×90=8640"},
    @{Row=10; Col=5; Text="85×51=4335"},

    @{Row=15; Col=1; Text="43×45=1935"},
    @{Row=15; Col=2; Text="38×57=2166"},
    @{Row=15; Col=3; Text="50×85=4250"},
    @{Row=15; Col=4; Text="92×54=4968"},
    @{Row=15; Col=5; Text="54×75=4050"},

    @{Row=20; Col=1; Text="41×46=1886"},
    @{Row=20; Col=2; Text="51×66=3366"},
    @{Row=20; Col=3; Text="34×85=2890"},
    @{Row=20; Col=4; Text="16×43=688"},
    @{Row=20; Col=5; Text="40×12=480"}
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
